$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$ReplaceAll = 2
$FindContinue = 1

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, $FindContinue, $false, $new, $ReplaceAll)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. "An order is related to a product..." customer -> customers
Replace-Text "An order is related to a product in a way that customer can place an order for a product in the gym." "An order is related to a product in a way that customers can place an order for a product in the gym."

# 2. Meals paragraph rewrite
Replace-Text ": There will be meal plans which will be recommended to the customer by the nutritionist which will be provided as a part of product entity. Meals do not need to be bound to nutritionist. A customer can also order meals as a part of product." ": There will be meal plans recommended to the customer by the nutritionist and provided as a part of the product entity. Meals do not need to be bound to a nutritionist. A customer can also order meals as a part of the product."

# 3. Product attributes list expanded
Replace-Text "A product can have attributes as product name, product price, etc." "A product can have attributes such as product name, product price, product type, product description  etc."

# 4. Recommends associative entity: which -> that
Replace-Text "Recommends is an associative entity which has recommended quantity and rec reason, date as attributes." "Recommends is an associative entity that has recommended quantity and rec reason, date as attributes."

# 5. Location attributes: add building size
Replace-Text "A location has state, area, city, street, zip code as attributes." "A location has state, area, city, street, zip code, building size as attributes."

# 6. Here purchases is/are
Replace-Text "Here purchases is an associative entity." "Here purchases are an associative entity."

# 7. Locker / equipment types -> weights and machines (leave the "ie" proofErr markers intact)
Replace-Text ", locker and equipments." ", locker, weights, and machines."

# 8. Customer can take services of a nutritionist
Replace-Text "A customer can take services of a nutritionist as well." "A customer can take the services of a nutritionist as well."

# 9. Employees relations: Manager -> Location, add "and"
Replace-Text "Employees have the relations with Department, Manager, Staff, Trainer, Nutritionist." "Employees have the relations with Department, Location, Staff, Trainer, and Nutritionist."

# 10. Customer relations: add recommends
Replace-Text "Customer has the relations with order, subscription, location, utilities, nutritionist and trainer." "Customer has the relations with order, subscription, recommends, location, utilities, nutritionist and trainer."

# 11. Remove the lastRenderedPageBreak artifact before "Customer, Location and Subscription Relationship:"
#     by re-writing the heading text in place (strips the stray rendering break element).
Replace-Text "Customer, Location and Subscription Relationship:" "Customer, Location and Subscription Relationship:"
